$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "About"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("About")

$ws1.Range("A1").Value = "WMITR Worker Marginal Income Tax Rate"
$ws1.Range("A3").Value = "Source:"
$ws1.Range("B3").Value = "Trading Economics"
$ws1.Range("B4").Value = 2020
$ws1.Range("B5").Value = "List of Countries by Personal Income Tax Rate "
$ws1.Range("B6").Value = "https://tradingeconomics.com/country-list/personal-income-tax-rate?continent=europe"
$ws1.Range("B7").Value = "European Union "
$ws1.Range("B8").Value = "Accessed 20th May 2020"

# Move the "Notes:" label from row 9 down to row 10 (copy formatting, then remove old row)
[void]$ws1.Range("A9").Copy($ws1.Range("A10"))
[void]$ws1.Rows(9).Clear()

$ws1.Range("A11").Value = "Note that this source does not specify whether this income tax applies to the average or workers, but it is likely that it is the average. "
$ws1.Range("A12").Value = "We assume that the overall difference between the overall average and worker-only value is insignificant."

# ---------------------------------------------------------------
# Sheet "WMITR"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("WMITR")

$ws2.Range("B2").Value = 0.386

# Update stored selections on each sheet, ending on "About" so it stays the active tab
[void]$ws2.Range("B3").Select()
[void]$ws1.Range("A13").Select()
